$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15944.639376
$ws.Range("D2").Value = 203.814385

$ws.Range("B3").Value = 2721.927865
$ws.Range("D3").Value = 17.396695
$ws.Range("E3").Value = 0

$ws.Range("B4").Value = 25972.750949
$ws.Range("C4").Value = 332

$ws.Range("G5").Value = -5.762183
$ws.Range("H5").Value = -8.968999
$ws.Range("I5").Value = -2.555368
$ws.Range("J5").Value = 0.000089

$ws.Range("G6").Value = -0.087215
$ws.Range("H6").Value = -3.455459
$ws.Range("I6").Value = 3.281028
$ws.Range("J6").Value = 0.997953

$ws.Range("G7").Value = 5.674968
$ws.Range("H7").Value = 3.142274
$ws.Range("I7").Value = 8.207662
$ws.Range("J7").Value = 0.000001
